# Edit: change the table style applied to the 4-column table on slide 6
# ("SOURCES OF FINANCE") from the deck's custom "Table_0" style to the
# built-in PowerPoint table style {F691944B-0757-4608-9A80-5585DB12D16C}.
#
# Table styles cannot be changed by assigning Table.Style directly (the
# host raises "Table styles cannot be assigned through a property - call
# Table.ApplyStyle(...) instead"), so we use Table.ApplyStyle with the
# style's GUID, exactly as PowerPoint's Table Styles gallery does when a
# user clicks a different style swatch.

$p = $ppt.ActivePresentation

# Locate the table shape on slide 6 robustly (rather than assuming a
# fixed shape index) by scanning for the shape that HasTable.
$slide = $p.Slides.Item(6)

$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

if ($tableShape -ne $null) {
    $table = $tableShape.Table
    $table.ApplyStyle("{F691944B-0757-4608-9A80-5585DB12D16C}")
}
